$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1 ("土地") -- rename headers, fix data typos, append new columns
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Extend the header row (B1:H1) with 7 new header cells (I1:O1), copying
# the existing header style/border along the way.
$ws1.Range("B1:H1").Copy($ws1.Range("I1"))

# Re-label the original headers (B1:H1)
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "area"
$ws1.Range("D1").Value = "share_portion"
$ws1.Range("E1").Value = "owner"
$ws1.Range("F1").Value = "register_date"
$ws1.Range("G1").Value = "register_reason"
$ws1.Range("H1").Value = "acquire_value"

# New headers (I1:O1)
$ws1.Range("I1").Value = "property_category"
$ws1.Range("J1").Value = "category"
$ws1.Range("K1").Value = "date"
$ws1.Range("L1").Value = "legislator_name"
$ws1.Range("M1").Value = "legislator_id"
$ws1.Range("N1").Value = "source_file"
$ws1.Range("O1").Value = "index"

# Clean up the data row's stray whitespace / punctuation typos (B2:H2)
$ws1.Range("B2").Value = "臺北市大安區龍泉段一小段02930000地號"
$ws1.Range("D2").Value = "100000分之16216"
$ws1.Range("E2").Value = "高金素梅"
$ws1.Range("F2").Value = "92年12月25日"
$ws1.Range("G2").Value = "買賣"
$ws1.Range("H2").Value = "25000000(土地建物與車位合併價）"

# New data cells (I2:O2)
$ws1.Range("I2").Value = "land"
$ws1.Range("J2").Value = "normal"

# K2 needs to hold the literal text "2011-11-22" -- a direct .Value
# assignment gets auto-coerced to a date serial by the smart setter, so
# round-trip it through a TEXT() formula and paste-special the computed
# value back in as plain text (keeps the shared-string literal and
# avoids minting a brand-new number-format style).
$ws1.Range("Q1").Formula = "=TEXT(""2011-11-22"",""@"")"
$ws1.Range("Q1").Copy()
$ws1.Range("K2").PasteSpecial(-4163)
$ws1.Range("Q1").ClearContents()

$ws1.Range("L2").Value = "高金素梅"
$ws1.Range("M2").Value = 926
$ws1.Range("N2").Value = "tmp2f3b1"
$ws1.Range("O2").Value = 14

# ----------------------------------------------------------------------
# Sheet 2 ("建物") -- same typo clean-up, shares strings with sheet 1
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "臺北市大安區龍泉段一小段05819000建號"
$ws2.Range("F2").Value = "92年12月25日"
$ws2.Range("H2").Value = "25000000(土地建物與車位合併價）"

# ----------------------------------------------------------------------
# Sheet 5 ("債務") -- strip stray spaces from addresses / date
# ----------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("D2").Value = "陳麗卿新北市泰山區明志路"
$ws5.Range("D3").Value = "石旭松新北市泰山區明志路"
$ws5.Range("F2").Value = "96年02月06日"
$ws5.Range("F3").Value = "96年02月06日"
